$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 4: swap the Assignment/Device Use (F/G) columns ---
# Before: F2="Button0"(s6)   G2="Feeder"(s3)
#         F4="Button0GND"(s11) G4=<empty>(s3)
# After:  F2="Button0GND"(s11)  G2=<absent>
#         F4="Button0"(s6)      G4="Feeder"(s3)
$ws.Range("F2:G2").Copy($ws.Range("Z1:AA1"))
$ws.Range("F4:G4").Copy($ws.Range("Z2:AA2"))

$ws.Range("Z1:AA1").Copy($ws.Range("F4:G4"))
$ws.Range("Z2").Copy($ws.Range("F2"))
$ws.Range("G2").Clear()

$ws.Range("Z1:AA2").Clear()

# --- Swap G17 <-> G25 (Switch_LeftRight_Left <-> Switch_LeftRight_Right) ---
$g17 = $ws.Range("G17").Value()
$g25 = $ws.Range("G25").Value()
$ws.Range("G17").Value = $g25
$ws.Range("G25").Value = $g17

# --- Swap G30 <-> G31 (Switch_Unassigned_1_OFF <-> Switch_Unassigned_1_ON) ---
$g30 = $ws.Range("G30").Value()
$g31 = $ws.Range("G31").Value()
$ws.Range("G30").Value = $g31
$ws.Range("G31").Value = $g30

# --- Swap G32 <-> G33 (Switch_Unassigned_2_OFF <-> Switch_Unassigned_2_ON) ---
$g32 = $ws.Range("G32").Value()
$g33 = $ws.Range("G33").Value()
$ws.Range("G32").Value = $g33
$ws.Range("G33").Value = $g32

# --- View state: zoom to 228% and move the selection from G34 to G32 ---
$excel.ActiveWindow.Zoom = 228
$ws.Range("G32").Select()
